# Apply the HotStock Top20 20250827 edits: update individual cell values
# in Sheet1 so that the resulting workbook matches the target revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value  = "拓维信息"
$ws.Range("C5").Value  = "寒武纪"
$ws.Range("C6").Value  = "成飞集成"
$ws.Range("C7").Value  = "利欧股份"
$ws.Range("A8").Value  = "英维克"
$ws.Range("A9").Value  = "合力泰"
$ws.Range("C11").Value = "工业富联"
$ws.Range("A12").Value = "剑桥科技"
$ws.Range("C12").Value = "吉视传媒"
$ws.Range("A13").Value = "吉视传媒"
$ws.Range("C13").Value = "鸿博股份"
$ws.Range("A14").Value = "新易盛"
$ws.Range("C14").Value = "东华软件"
$ws.Range("A15").Value = "中油资本"
$ws.Range("C16").Value = "新易盛"
$ws.Range("C17").Value = "天融信"
$ws.Range("C18").Value = "步步高"
$ws.Range("A19").Value = "成飞集成"
$ws.Range("C19").Value = "润和软件"
$ws.Range("A20").Value = "奋达科技"
$ws.Range("C21").Value = "卓翼科技"
